$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17229425907135
$ws.Range("B1").Value = 2.438127040863037
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.366897106170654
$ws.Range("E1").Value = 1.235362648963928
